$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G - "Recorded By"
    $val = $cell.Value()
    if ($null -ne $val -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Length -gt 1) {
            $rotated = @($parts[$parts.Length - 1]) + $parts[0..($parts.Length - 2)]
            $newVal = [string]::Join(", ", $rotated)
            $cell.Value = $newVal
        }
    }
}
